# The canonical edit swaps the raw contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml: theme1.xml (previously the default "Office Theme"
# colour scheme) becomes the "Integral" colour scheme, and theme2.xml
# (previously "Integral", the scheme actually driving the slide master /
# slides) becomes the default "Office Theme" colour scheme. Font scheme and
# format scheme are identical between the two theme parts already, so the
# only observable difference is the 10 theme colours that differ from
# "Integral" (dk1/lt1 are black/white in both, so they are unaffected).
#
# This COM host models a single live theme (the one driving the slide
# master/slides), reachable through ThemeColorScheme off the design /
# master / slide. Re-pointing that theme's colours to the "Office Theme"
# palette reproduces the half of the swap that is visible through the
# object model.

function ToComColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# ThemeColorScheme item order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6,
# 11 hlink, 12 folHlink.
$officeTheme = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

foreach ($i in 1..12) {
    $tcs.Item($i).RGB = ToComColor $officeTheme[$i]
}
